$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2002047.8
$ws.Range("I40").Value = 2002047.8
$ws.Range("K40").Value = 2002047.8
$ws.Range("M40").Value = -2001872.8
$ws.Range("H48").Value = 5500
$ws.Range("I48").Value = 5000
$ws.Range("J48").Value = 6000
$ws.Range("K48").Value = 15000
$ws.Range("L48").Value = 18000
$ws.Range("M48").Value = -14708
$ws.Range("N48").Value = -18584
$ws.Range("H56").Value = 5500
$ws.Range("I56").Value = 5000
$ws.Range("J56").Value = 6000
$ws.Range("K56").Value = 15000
$ws.Range("L56").Value = 18000
$ws.Range("M56").Value = -14466
$ws.Range("N56").Value = -19068
$ws.Range("H62").Value = 5240.1177
$ws.Range("I62").Value = 3266.4285
$ws.Range("J62").Value = 6621.7
$ws.Range("K62").Value = 3266.4285
$ws.Range("L62").Value = 6621.7
$ws.Range("M62").Value = -2642.4285
$ws.Range("N62").Value = -7869.7
$ws.Range("H65").Value = 5240.1177
$ws.Range("I65").Value = 3266.4285
$ws.Range("J65").Value = 6621.7
$ws.Range("K65").Value = 16332.1425
$ws.Range("L65").Value = 33108.5
$ws.Range("M65").Value = -13212.1425
$ws.Range("N65").Value = -39348.5
$ws.Range("H107").Value = 1365.25
$ws.Range("I107").Value = 1244.5625
$ws.Range("J107").Value = 1848
$ws.Range("K107").Value = 1244.5625
$ws.Range("L107").Value = 1848
$ws.Range("M107").Value = 675.4375
$ws.Range("N107").Value = -5688
$ws.Range("H129").Value = 2557.6667
$ws.Range("I129").Value = 1086.75
$ws.Range("J129").Value = 5499.5
$ws.Range("K129").Value = 3260.25
$ws.Range("L129").Value = 16498.5
$ws.Range("M129").Value = 1739.75
$ws.Range("N129").Value = -26498.5
$ws.Range("H132").Value = 2129.125
$ws.Range("I132").Value = 1725.6666
$ws.Range("K132").Value = 5176.9998
$ws.Range("M132").Value = -2646.9998
$ws.Range("H141").Value = 1590.8422
$ws.Range("I141").Value = 1590.8422
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 4772.5266
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 407.4733999999999
$ws.Range("N141").ClearContents()
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6285.8184
$ws.Range("I32").Value = 4576.4644
$ws.Range("K32").Value = 4576.4644
$ws.Range("M32").Value = -4289.4644
$ws.Range("H63").Value = 6153.1055
$ws.Range("I63").Value = 4380.8
$ws.Range("J63").Value = 6786.0713
$ws.Range("K63").Value = 4380.8
$ws.Range("L63").Value = 6786.0713
$ws.Range("M63").Value = -3694.8
$ws.Range("N63").Value = -8158.0713
$ws.Range("H66").Value = 6153.1055
$ws.Range("I66").Value = 4380.8
$ws.Range("J66").Value = 6786.0713
$ws.Range("K66").Value = 21904
$ws.Range("L66").Value = 33930.35649999999
$ws.Range("M66").Value = -18472
$ws.Range("N66").Value = -40794.35649999999
$ws.Range("H122").Value = 4298.3076
$ws.Range("I122").Value = 3500
$ws.Range("K122").Value = 10500
$ws.Range("M122").Value = -8050
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1348.4286
$ws.Range("I20").Value = 1278.8182
$ws.Range("K20").Value = 1278.8182
$ws.Range("M20").Value = -1031.8182
$ws.Range("H42").Value = 200000
$ws.Range("J42").Value = 200000
$ws.Range("L42").Value = 200000
$ws.Range("N42").Value = -200656
$ws.Range("H86").Value = 18111.334
$ws.Range("I86").Value = 9617.75
$ws.Range("K86").Value = 9617.75
$ws.Range("M86").Value = -8494.75
$ws.Range("H89").Value = 18111.334
$ws.Range("I89").Value = 9617.75
$ws.Range("K89").Value = 48088.75
$ws.Range("M89").Value = -42472.75
$ws.Range("H94").Value = 1218.5
$ws.Range("I94").Value = 812.6667
$ws.Range("K94").Value = 812.6667
$ws.Range("M94").Value = -361.6667
$ws.Range("H134").Value = 4439
$ws.Range("I134").Value = 4392.25
$ws.Range("K134").Value = 13176.75
$ws.Range("M134").Value = -10641.75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 1649.5
$ws.Range("I41").Value = 799.5
$ws.Range("J41").Value = 2499.5
$ws.Range("K41").Value = 2398.5
$ws.Range("L41").Value = 7498.5
$ws.Range("M41").Value = -2060.5
$ws.Range("N41").Value = -8174.5
$ws.Range("H87").Value = 1499
$ws.Range("I87").Value = 1499
$ws.Range("K87").Value = 4497
$ws.Range("M87").Value = -3249
$ws.Range("H90").Value = 1499
$ws.Range("I90").Value = 1499
$ws.Range("K90").Value = 13491
$ws.Range("M90").Value = -7251
$ws.Range("H107").Value = 690
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 690
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 2070
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -5910
$ws.Range("H129").Value = 3942.8276
$ws.Range("I129").Value = 4588.625
$ws.Range("J129").Value = 3696.8096
$ws.Range("K129").Value = 13765.875
$ws.Range("L129").Value = 11090.4288
$ws.Range("M129").Value = -8765.875
$ws.Range("N129").Value = -21090.4288
$ws.Range("H131").Value = 32943.86
$ws.Range("I131").Value = 116575.445
$ws.Range("J131").Value = 5066.6665
$ws.Range("K131").Value = 349726.335
$ws.Range("L131").Value = 15199.9995
$ws.Range("M131").Value = -344686.335
$ws.Range("N131").Value = -25279.9995
$ws.Range("H132").Value = 1583.6364
$ws.Range("J132").Value = 1505.1666
$ws.Range("L132").Value = 13546.4994
$ws.Range("N132").Value = -18606.4994
$ws.Range("H138").Value = 2304.1667
$ws.Range("J138").Value = 2998.5
$ws.Range("L138").Value = 8995.5
$ws.Range("N138").Value = -19275.5
$ws.Range("H139").Value = 4086.3333
$ws.Range("J139").Value = 3333
$ws.Range("L139").Value = 9999
$ws.Range("N139").Value = -20279
$ws.Range("H140").Value = 999
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3397.5881
$ws.Range("I97").Value = 3350.6365
$ws.Range("K97").Value = 3350.6365
$ws.Range("M97").Value = -2854.6365
$ws.Range("H102").Value = 2427.1155
$ws.Range("I102").Value = 1785.7
$ws.Range("K102").Value = 1785.7
$ws.Range("M102").Value = -163.7
$ws.Range("H122").Value = 3460.5
$ws.Range("I122").Value = 2295.111
$ws.Range("K122").Value = 6885.333
$ws.Range("M122").Value = -4435.333
$ws.Range("H126").Value = 4590.5557
$ws.Range("I126").Value = 4367.4165
$ws.Range("K126").Value = 13102.2495
$ws.Range("M126").Value = -10632.2495
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3628.7058
$ws.Range("I7").Value = 3388.4443
$ws.Range("J7").Value = 3899
$ws.Range("K7").Value = 3388.4443
$ws.Range("L7").Value = 3899
$ws.Range("M7").Value = -3276.4443
$ws.Range("N7").Value = -4123
$ws.Range("H40").Value = 3303.4546
$ws.Range("I40").Value = 3409.842
$ws.Range("J40").Value = 2629.6667
$ws.Range("K40").Value = 3409.842
$ws.Range("L40").Value = 2629.6667
$ws.Range("M40").Value = -3273.842
$ws.Range("N40").Value = -2901.6667
$ws.Range("H61").Value = 3751.4075
$ws.Range("I61").Value = 3330.5625
$ws.Range("J61").Value = 4363.5454
$ws.Range("K61").Value = 3330.5625
$ws.Range("L61").Value = 4363.5454
$ws.Range("M61").Value = -3128.5625
$ws.Range("N61").Value = -4767.5454
$ws.Range("H75").Value = 58000
$ws.Range("I75").Value = 57000
$ws.Range("J75").Value = 60000
$ws.Range("K75").Value = 57000
$ws.Range("L75").Value = 60000
$ws.Range("M75").Value = -56064
$ws.Range("N75").Value = -61872
$ws.Range("H78").Value = 58000
$ws.Range("I78").Value = 57000
$ws.Range("J78").Value = 60000
$ws.Range("K78").Value = 171000
$ws.Range("L78").Value = 180000
$ws.Range("M78").Value = -166320
$ws.Range("N78").Value = -189360
$ws.Range("H106").Value = 167814.33
$ws.Range("J106").Value = 167814.33
$ws.Range("L106").Value = 167814.33
$ws.Range("N106").Value = -170338.33
$ws.Range("H113").Value = 3751.4075
$ws.Range("I113").Value = 3330.5625
$ws.Range("J113").Value = 4363.5454
$ws.Range("K113").Value = 3330.5625
$ws.Range("L113").Value = 4363.5454
$ws.Range("M113").Value = -1160.5625
$ws.Range("N113").Value = -8703.545399999999
$ws.Range("H122").Value = 4342.926
$ws.Range("I122").Value = 3355.375
$ws.Range("J122").Value = 4758.737
$ws.Range("K122").Value = 10066.125
$ws.Range("L122").Value = 14276.211
$ws.Range("M122").Value = -7616.125
$ws.Range("N122").Value = -19176.211
$ws.Range("H126").Value = 3628.7058
$ws.Range("I126").Value = 3388.4443
$ws.Range("J126").Value = 3899
$ws.Range("K126").Value = 10165.3329
$ws.Range("L126").Value = 11697
$ws.Range("M126").Value = -7695.332900000001
$ws.Range("N126").Value = -16637
$ws.Range("H136").Value = 4888.5
$ws.Range("I136").Value = 1875
$ws.Range("J136").Value = 7299.3
$ws.Range("K136").Value = 5625
$ws.Range("L136").Value = 21897.9
$ws.Range("M136").Value = -3075
$ws.Range("N136").Value = -26997.9
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3400.25
$ws.Range("I122").Value = 2867
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 8601
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -6151
$ws.Range("N122").Value = -19900
